$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (Subject) from 18 to 22 characters
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668

# Update Subject labels and Recorded-By email orderings per daily attendance refresh
$ws.Range("C7").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C8").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G9").Value = 'wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("C21").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C22").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G29").Value = 'Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("C38").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C39").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G40").Value = 'wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Range("C52").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C53").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G60").Value = 'Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G64").Value = 'mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Range("G65").Value = 'Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("C69").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C70").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G71").Value = 'Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("C83").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C84").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G91").Value = 'Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("C100").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C101").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G102").Value = 'Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("C114").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C115").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G122").Value = 'Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("C131").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C132").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G133").Value = 'Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("C145").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C146").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G153").Value = 'Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G157").Value = 'servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G158").Value = 'Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range("C162").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("C163").Value = 'BIOCHEMISTRY LAB/CBL'
$ws.Range("G164").Value = 'Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("C176").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("C177").Value = 'PARASITOLOGY SGD/POS'
$ws.Range("G184").Value = 'maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
